$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new Column K ("Is Active?") ------------------------------------
# Header K1: copy format from J1 (bold/underlined header style) then set text
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K1").Value = "Is Active?"

# K2:K7 - copy boolean cell format from L2 (fill style) and set to TRUE
$boolRows = 2, 3, 5, 6, 7
foreach ($r in $boolRows) {
    $ws.Range("L$r").Copy() | Out-Null
    $ws.Range("K$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("K$r").Value = $true
}

# --- Add new data rows 8 and 9 -------------------------------------------
# Row 8: ToggleToActiveGroup
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value = "ToggleToActiveGroup"

$ws.Range("H7").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null
$ws.Range("H8").Value = 7

$ws.Range("J7").Copy() | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").Value = "ToggleToActiveGroup"

$ws.Range("L7").Copy() | Out-Null
$ws.Range("K8").PasteSpecial(-4122) | Out-Null
$ws.Range("K8").Value = $true

$ws.Range("L7").Copy() | Out-Null
$ws.Range("L8").PasteSpecial(-4122) | Out-Null
$ws.Range("L8").Value = $false

# Row 9: ToggleToInactiveGroup
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Value = "ToggleToInactiveGroup"

$ws.Range("H7").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("H9").Value = 8

$ws.Range("J7").Copy() | Out-Null
$ws.Range("J9").PasteSpecial(-4122) | Out-Null
$ws.Range("J9").Value = "ToggleToInactiveGroup"

$ws.Range("L7").Copy() | Out-Null
$ws.Range("K9").PasteSpecial(-4122) | Out-Null
$ws.Range("K9").Value = $false

$ws.Range("L7").Copy() | Out-Null
$ws.Range("L9").PasteSpecial(-4122) | Out-Null
$ws.Range("L9").Value = $false

# --- Update sheet view: move the selection to J9 ------------------------
$ws.Range("J9").Select() | Out-Null
